# Handback status report: add a new handed-back file
# (4f5e71f0-3a6a-4d21-ac63-8bc858d9d78c.md) to the Overview, zh-cn and
# de-de tables, mirroring the two rows already present in each.
#
# NOTE: "True" / "False" / "" need a leading apostrophe so Excel stores
# them as plain text (shared-string) cells -- exactly like the existing
# rows -- instead of auto-converting them to boolean cells.

$wb = $excel.ActiveWorkbook

$fileId = "4f5e71f0-3a6a-4d21-ac63-8bc858d9d78c"
$mdName = "$fileId.md"
$mdPath = "e2e\$fileId.md"
$zhXlf  = "$fileId.05f29578d6a440f4d81d5703299c6f38047b1100.zh-cn.xlf"
$deXlf  = "$fileId.05f29578d6a440f4d81d5703299c6f38047b1100.de-de.xlf"

$mdUrlOverview = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b8b9f94453e6f2fac80f79f7be1afc5081a86a6/e2e/$mdName"
$mdUrlZhMain   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b8b9f94453e6f2fac80f79f7be1afc5081a86a6/e2e/$mdName"
$mdUrlDeMain   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b8b9f94453e6f2fac80f79f7be1afc5081a86a6/e2e/$mdName"
$mdUrlZh       = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/78087c533e0802bfe192977e3cf9157b5502fe51/e2e/$mdName"
$mdUrlDe       = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0d97826f71631239c6fea01b2af311c0fc40acf9/e2e/$mdName"

$statusSync = "Handed back: in sync with en-US"
$trueText   = "'True"
$falseText  = "'False"
$emptyText  = "'"

# ---------------------------------------------------------------------
# Sheet "Overview" -> table "Overview" (A1:G3 -> A1:G4)
# ---------------------------------------------------------------------
$wsOverview  = $wb.Worksheets.Item("Overview")
$loOverview  = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()
$r = 1 + $rowOverview.Index

$wsOverview.Cells.Item($r, 1).Value = $mdName
$wsOverview.Cells.Item($r, 3).Value = ".md"
$wsOverview.Cells.Item($r, 5).Value = $statusSync
$wsOverview.Cells.Item($r, 6).Value = $statusSync
$wsOverview.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Cells.Item($r, 7).Value = "2016-09-02 20:51:47"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($r, 2), $mdUrlOverview, $null, $null, $mdPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> table "zh-cn" (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------
$wsZh  = $wb.Worksheets.Item("zh-cn")
$loZh  = $wsZh.ListObjects.Item("zh-cn")
$rowZh = $loZh.ListRows.Add()
$r = 1 + $rowZh.Index

$wsZh.Cells.Item($r, 2).Value  = ".md"
$wsZh.Cells.Item($r, 3).Value  = $statusSync
$wsZh.Cells.Item($r, 4).Value  = "e2e"
$wsZh.Cells.Item($r, 5).Value  = "ht"
$wsZh.Cells.Item($r, 6).Value  = $trueText
$wsZh.Cells.Item($r, 7).Value  = $zhXlf
$wsZh.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($r, 8).Value  = "2016-09-02 20:51:41"
$wsZh.Cells.Item($r, 10).Value = $zhXlf
$wsZh.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($r, 11).Value = "2016-09-02 20:51:59"
$wsZh.Cells.Item($r, 12).Value = $emptyText
$wsZh.Cells.Item($r, 13).Value = $trueText
$wsZh.Cells.Item($r, 14).Value = $emptyText
$wsZh.Cells.Item($r, 15).Value = $falseText
$wsZh.Cells.Item($r, 16).Value = $emptyText

$wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 1), $mdUrlZhMain, $null, $null, $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 9), $mdUrlZh, $null, $null, $mdName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" -> table "de-de" (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------
$wsDe  = $wb.Worksheets.Item("de-de")
$loDe  = $wsDe.ListObjects.Item("de-de")
$rowDe = $loDe.ListRows.Add()
$r = 1 + $rowDe.Index

$wsDe.Cells.Item($r, 2).Value  = ".md"
$wsDe.Cells.Item($r, 3).Value  = $statusSync
$wsDe.Cells.Item($r, 4).Value  = "e2e"
$wsDe.Cells.Item($r, 5).Value  = "ht"
$wsDe.Cells.Item($r, 6).Value  = $trueText
$wsDe.Cells.Item($r, 7).Value  = $deXlf
$wsDe.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($r, 8).Value  = "2016-09-02 20:51:47"
$wsDe.Cells.Item($r, 10).Value = $deXlf
$wsDe.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($r, 11).Value = "2016-09-02 20:52:15"
$wsDe.Cells.Item($r, 12).Value = $emptyText
$wsDe.Cells.Item($r, 13).Value = $trueText
$wsDe.Cells.Item($r, 14).Value = $emptyText
$wsDe.Cells.Item($r, 15).Value = $falseText
$wsDe.Cells.Item($r, 16).Value = $emptyText

$wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 1), $mdUrlDeMain, $null, $null, $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 9), $mdUrlDe, $null, $null, $mdName) | Out-Null
